# Weekly data refresh: insert the newest week's record at the top of the
# data block (row 4, right after the two most-recent existing rows) and
# push every older record down by one row. The previously-oldest row
# (old row 63) ends up as the new last row (row 64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4..63 down to 5..64, creating a blank row 4 for the new record.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A4").Value2 = 6
$ws.Range("B4").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value2 = "Metropolitana"
$ws.Range("D4").Value2 = 44699
$ws.Range("E4").Value2 = 13
$ws.Range("F4").Value2 = "Fruta"
$ws.Range("G4").Value2 = 100108
$ws.Range("H4").Value2 = "Tropicales y subtropicales"
$ws.Range("I4").Value2 = 100108007
$ws.Range("J4").Value2 = "Coco"
$ws.Range("K4").Value2 = "Sin especificar"
$ws.Range("L4").Value2 = "Primera"
$ws.Range("M4").Value2 = 500
$ws.Range("N4").Value2 = 23000
$ws.Range("O4").Value2 = 24000
$ws.Range("P4").Value2 = 23500
$ws.Range("Q4").Value2 = "$/malla 20 unidades"
$ws.Range("R4").Value2 = "Perú"
$ws.Range("S4").Value2 = 1175
$ws.Range("T4").Value2 = 20
